$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.181.43"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.423.22"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.27"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.53"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.06"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "2.853.93"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "60.086.62"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "2.414.84"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.43"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.67"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.71"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.40"
$ws.Range("E27").Value = "  +3.38%  "
$ws.Range("D28").Value = "0.0₃0776"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.02"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.405"
$ws.Range("E33").Value = "  -4.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.61"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "333.56"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.89"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.52"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.66"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.15"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0224"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("E51").Value = "  -1.17%  "
